# Update profit.py after running on 2025-08-29
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Append the new day's profit figure to Sheet1 (A1:B11 -> A1:B12)
# (force Text so the date-like string isn't auto-coerced into a date
#  serial, then restore the default style so no residual formatting
#  is left behind on the cell)
$ws1.Range("A12").NumberFormat = "@"
$ws1.Range("A12").Value = "08/29/2025"
$ws1.Range("A12").Style = "Normal"
$ws1.Range("B12").Value = 11839.93

# Refresh the rolling single-row stats sheet for the new date
$ws2.Range("A1").NumberFormat = "@"
$ws2.Range("A1").Value = "08/29/2025"
$ws2.Range("A1").Style = "Normal"
$ws2.Range("B1").Value = 0.1104204078508326
$ws2.Range("C1").Value = 0.8895795921491674
